# Update the header date, then update each division problem in the table.
# Replacements are applied in document order (top-to-bottom, left-to-right)
# so that a value introduced by an earlier replacement (e.g. "84÷6=" being
# written by the 83÷4= -> 84÷6= substitution) is never accidentally matched
# by a later Find call looking for that same original text elsewhere.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-19 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-20 Thursday", 2) | Out-Null
$d.Content.Find.Execute("50÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷9=", 2) | Out-Null
$d.Content.Find.Execute("16÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷2=", 2) | Out-Null
$d.Content.Find.Execute("41÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=", 2) | Out-Null
$d.Content.Find.Execute("64÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷7=", 2) | Out-Null
$d.Content.Find.Execute("84÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷9=", 2) | Out-Null
$d.Content.Find.Execute("17÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=", 2) | Out-Null
$d.Content.Find.Execute("33÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷9=", 2) | Out-Null
$d.Content.Find.Execute("46÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=", 2) | Out-Null
$d.Content.Find.Execute("89÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=", 2) | Out-Null
$d.Content.Find.Execute("66÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷8=", 2) | Out-Null
$d.Content.Find.Execute("24÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷8=", 2) | Out-Null
$d.Content.Find.Execute("67÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷9=", 2) | Out-Null
$d.Content.Find.Execute("83÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷6=", 2) | Out-Null
$d.Content.Find.Execute("31÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷4=", 2) | Out-Null
$d.Content.Find.Execute("48÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷3=", 2) | Out-Null
$d.Content.Find.Execute("78÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=", 2) | Out-Null
$d.Content.Find.Execute("26÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷8=", 2) | Out-Null
$d.Content.Find.Execute("74÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷3=", 2) | Out-Null
$d.Content.Find.Execute("47÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=", 2) | Out-Null
$d.Content.Find.Execute("46÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=", 2) | Out-Null
$d.Content.Find.Execute("46÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷2=", 2) | Out-Null
$d.Content.Find.Execute("15÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷7=", 2) | Out-Null
$d.Content.Find.Execute("80÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=", 2) | Out-Null
$d.Content.Find.Execute("27÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷3=", 2) | Out-Null
$d.Content.Find.Execute("99÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=", 2) | Out-Null
